$d = $word.ActiveDocument

$shape = $d.Shapes.Item(1)
Write-Output ("Before Left: " + $shape.Left)
$shape.Left = -34.866692913385826
Write-Output ("After Left1: " + $shape.Left)
$shape.Top = -16.8
Write-Output ("After Top1: " + $shape.Top)

$table = $d.Tables.Item(1)
$table.Rows.VerticalPosition = 1708
Write-Output ("After VerticalPosition: " + $table.Rows.VerticalPosition)

Write-Output ("Final Left: " + $shape.Left)
Write-Output ("Final Top: " + $shape.Top)
